# Ajout amélioration personnage + arme et début equipement
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StatsPersonnages")

# weaponName (column J) values for rows 3,5,7,9,11,13,15,17
$ws.Range("J3").Value = "Hache"
$ws.Range("J5").Value = "Dague"
$ws.Range("J7").Value = "BouleElectrique"
$ws.Range("J9").Value = "Epee"
$ws.Range("J11").Value = "Flechette"
$ws.Range("J13").Value = "Lance"
$ws.Range("J15").Value = "BouleDeFeu"
$ws.Range("J17").Value = "Shuriken"

# name (column O) values for rows 3,5,7,9,11,13,15,17
$ws.Range("O3").Value = "Djo"
$ws.Range("O5").Value = "Kral"
$ws.Range("O7").Value = "Cork"
$ws.Range("O9").Value = "Nato"
$ws.Range("O11").Value = "Cyrdin"
$ws.Range("O13").Value = "Galdir"
$ws.Range("O15").Value = "Swift"
$ws.Range("O17").Value = "Shoginawa"

# Update the active cell selection to match the recorded end state
$ws.Range("G14").Select()
